# Small corrections in mapping of land regionalized flows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = "AE"
$ws.Range("D3").Value = "AE"
$ws.Range("C5").Value = "AG"
$ws.Range("C17").Value = "BF"
$ws.Range("C20").Value = "BI"
$ws.Range("C59").Value = "BZ"
$ws.Range("C75").Value = "CF"
$ws.Range("C124").Value = "CW"
$ws.Range("C125").Value = "CY"
$ws.Range("C134").Value = "EG"
$ws.Range("C144").Value = "FJ"
$ws.Range("C152").Value = "GM"
$ws.Range("C155").Value = "GY"
$ws.Range("B156").Value = "HK"
$ws.Range("C170").Value = "IL"
$ws.Range("C212").Value = "IS"
$ws.Range("C215").Value = "JO"
$ws.Range("C222").Value = "KW"
$ws.Range("C224").Value = "LA"
$ws.Range("C225").Value = "LB"
$ws.Range("C226").Value = "LC"
$ws.Range("C228").Value = "LS"
$ws.Range("C232").Value = "LY"
$ws.Range("C238").Value = "ML"
$ws.Range("C241").Value = "MR"
$ws.Range("C242").Value = "MS"
$ws.Range("C243").Value = "MT"
$ws.Range("C246").Value = "MW"
$ws.Range("C260").Value = "OM"
$ws.Range("C270").Value = "QA"
$ws.Range("C286").Value = "SA"
$ws.Range("C288").Value = "SC"
$ws.Range("C295").Value = "SR"
$ws.Range("C298").Value = "SY"
$ws.Range("C299").Value = "SZ"
$ws.Range("C302").Value = "TJ"
$ws.Range("C303").Value = "TL"
$ws.Range("C313").Value = "UG"
$ws.Range("C327").Value = "UZ"
$ws.Range("C333").Value = "WS"
$ws.Range("D333").Value = "RAS"
$ws.Range("E333").Value = "RAS"
$ws.Range("C335").Value = "YE"

[void]$ws.Range("A1").Select()
